# Update variable description metadata (header label + a few cell fixes)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header column F: "Note" -> "Description"
$ws.Range("F1").Value = "Description"

# Example column: turn text dates into real dates (keep yyyy-mm-dd display format)
$ws.Range("E2").Value = (Get-Date -Year 2000 -Month 3 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E5").Value = (Get-Date -Year 2000 -Month 3 -Day 2 -Hour 0 -Minute 0 -Second 0)

# ShortWave / LongWave units: drop the superscript formatting on "W/m2"
$ws.Range("C6").Value = "W/m2"
$ws.Range("D6").Value = "decimal"
$ws.Range("C7").Value = "W/m2"
$ws.Range("D7").Value = "decimal"

# RelHum: add missing Format value
$ws.Range("D9").Value = "decimal"

# WindSpeed: remove stray leading space, add missing Format value
$ws.Range("B10").Value = "WindSpeed"
$ws.Range("D10").Value = "decimal $([char]0x2265) 0"

# Column F width to fit the new "Description" column content
$ws.Columns.Item(6).ColumnWidth = 41

# Selection now spans the whole used range instead of a single cell
$ws.Range("A1:F12").Select()
